$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info ---
# Project name
$ws.Range("B2").Value = "Sneeze Partition Installation"
# Date (stored as literal text in the template, not a real date - force
# text format so COM doesn't coerce the date-looking string into a serial)
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "05/03/2021"
# Contract No.
$ws.Range("B3").Value = "9055-017.00.01"
# Project ID
$ws.Range("F3").Value = 6261
# Start time
$ws.Range("F4").Value = "05:00"
# Weather
$ws.Range("B5").Value = "Sunny"
# End time
$ws.Range("F5").Value = "16:00"

# --- Contractor / Trade / Manpower / Equipment / Work performed table ---
# Row 8
$ws.Range("A8").Value = "Exbon Development Inc."
$ws.Range("B8").Value = "Carpenter"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = "Forklift"
$ws.Range("F8").Value = "Sneeze Partition Installation"

# Row 9
$ws.Range("A9").Value = "JPUS"
$ws.Range("B9").Value = "Laborer"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = "Dump Truck"
$ws.Range("F9").Value = "Existing Partition Removal"

# --- Tests & Inspections / Correctional Items / Note free-text sections ---
$ws.Range("A18").Value = "Inspection is scheduled on May 4, 2021 at 3PM. "
$ws.Range("A23").Value = "New frosted panel has a crack. Needs to be replaced."
$ws.Range("A28").Value = "All punchwork need to be performed during off-hours."
